# Apply the updated cryptocurrency market data (prices / 1h volume %, and the
# reordering of the OKB/Dogecoin and Frax/TrustWalletToken rows) to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = "28.405.11"   # D2: "28.407.70" -> "28.405.11"
$ws.Cells.Item(2, 5).Value2 = "  +6.03%  "   # E2: "  +6.02%  " -> "  +6.03%  "

# Row 3
$ws.Cells.Item(3, 4).Value2 = "1.815.59"   # D3: "1.816.55" -> "1.815.59"
$ws.Cells.Item(3, 5).Value2 = "  +5.49%  "   # E3: "  +5.50%  " -> "  +5.49%  "

# Row 4
$c = $ws.Cells.Item(4, 4)   # D4: "0.9998" -> "1.000"
$c.NumberFormat = "@"
$c.Value2 = "1.000"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value2 = "  -0.25%  "   # E4: "  -0.33%  " -> "  -0.25%  "

# Row 5
$c = $ws.Cells.Item(5, 4)   # D5: "318.05" -> "318.07"
$c.NumberFormat = "@"
$c.Value2 = "318.07"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = "  +2.72%  "   # E5: "  +2.73%  " -> "  +2.72%  "

# Row 6
$ws.Cells.Item(6, 5).Value2 = "  -0.26%  "   # E6: "  -0.29%  " -> "  -0.26%  "

# Row 7
$c = $ws.Cells.Item(7, 4)   # D7: "0.5703" -> "0.5706"
$c.NumberFormat = "@"
$c.Value2 = "0.5706"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value2 = "  +17.49%  "   # E7: "  +17.65%  " -> "  +17.49%  "

# Row 8
$c = $ws.Cells.Item(8, 4)   # D8: "0.3869" -> "0.3871"
$c.NumberFormat = "@"
$c.Value2 = "0.3871"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value2 = "  +11.35%  "   # E8: "  +11.34%  " -> "  +11.35%  "

# Row 9
$ws.Cells.Item(9, 2).Value2 = "Dogecoin"   # B9: "OKB" -> "Dogecoin"
$ws.Cells.Item(9, 3).Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"   # C9: "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" -> "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Cells.Item(9, 4)   # D9: "43.24" -> "0.07621"
$c.NumberFormat = "@"
$c.Value2 = "0.07621"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value2 = "  +5.58%  "   # E9: "  +1.25%  " -> "  +5.58%  "

# Row 10
$ws.Cells.Item(10, 2).Value2 = "OKB"   # B10: "Dogecoin" -> "OKB"
$ws.Cells.Item(10, 3).Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"   # C10: "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge" -> "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Cells.Item(10, 4)   # D10: "0.07619" -> "43.16"
$c.NumberFormat = "@"
$c.Value2 = "43.16"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value2 = "  +1.06%  "   # E10: "  +5.49%  " -> "  +1.06%  "

# Row 11
$c = $ws.Cells.Item(11, 4)   # D11: "1.139" -> "1.140"
$c.NumberFormat = "@"
$c.Value2 = "1.140"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = "  +8.75%  "   # E11: "  +8.74%  " -> "  +8.75%  "

# Row 12
$c = $ws.Cells.Item(12, 4)   # D12: "21.28" -> "21.29"
$c.NumberFormat = "@"
$c.Value2 = "21.29"
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 4)   # D13: "0.9998" -> "0.9996"
$c.NumberFormat = "@"
$c.Value2 = "0.9996"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value2 = "  -0.29%  "   # E13: "  -0.31%  " -> "  -0.29%  "

# Row 14
$c = $ws.Cells.Item(14, 4)   # D14: "6.273" -> "6.269"
$c.NumberFormat = "@"
$c.Value2 = "6.269"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value2 = "  +7.04%  "   # E14: "  +7.12%  " -> "  +7.04%  "

# Row 15
$ws.Cells.Item(15, 4).Value2 = "1.813.64"   # D15: "1.809.45" -> "1.813.64"
$ws.Cells.Item(15, 5).Value2 = "  +5.33%  "   # E15: "  +5.47%  " -> "  +5.33%  "

# Row 16
$c = $ws.Cells.Item(16, 4)   # D16: "7.294" -> "7.284"
$c.NumberFormat = "@"
$c.Value2 = "7.284"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value2 = "  +7.11%  "   # E16: "  +7.21%  " -> "  +7.11%  "

# Row 17
$c = $ws.Cells.Item(17, 4)   # D17: "92.21" -> "92.15"
$c.NumberFormat = "@"
$c.Value2 = "92.15"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value2 = "  +6.75%  "   # E17: "  +6.77%  " -> "  +6.75%  "

# Row 18
$c = $ws.Cells.Item(18, 4)   # D18: "0.00001077" -> "0.00001078"
$c.NumberFormat = "@"
$c.Value2 = "0.00001078"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value2 = "  +4.06%  "   # E18: "  +3.84%  " -> "  +4.06%  "

# Row 19
$c = $ws.Cells.Item(19, 4)   # D19: "0.06471" -> "0.06476"
$c.NumberFormat = "@"
$c.Value2 = "0.06476"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value2 = "  +1.06%  "   # E19: "  +1.04%  " -> "  +1.06%  "

# Row 20
$c = $ws.Cells.Item(20, 4)   # D20: "0.9995" -> "0.9996"
$c.NumberFormat = "@"
$c.Value2 = "0.9996"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = "  -0.26%  "   # E20: "  -0.25%  " -> "  -0.26%  "

# Row 21
$c = $ws.Cells.Item(21, 4)   # D21: "17.34" -> "17.36"
$c.NumberFormat = "@"
$c.Value2 = "17.36"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value2 = "  +5.11%  "   # E21: "  +4.96%  " -> "  +5.11%  "

# Row 22
$c = $ws.Cells.Item(22, 4)   # D22: "6.006" -> "6.003"
$c.NumberFormat = "@"
$c.Value2 = "6.003"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value2 = "  +5.24%  "   # E22: "  +5.20%  " -> "  +5.24%  "

# Row 23
$ws.Cells.Item(23, 4).Value2 = "28.421.87"   # D23: "28.418.01" -> "28.421.87"
$ws.Cells.Item(23, 5).Value2 = "  +5.81%  "   # E23: "  +5.71%  " -> "  +5.81%  "

# Row 24
$ws.Cells.Item(24, 5).Value2 = "  +3.75%  "   # E24: "  +3.71%  " -> "  +3.75%  "

# Row 25
$c = $ws.Cells.Item(25, 4)   # D25: "2.123" -> "2.127"
$c.NumberFormat = "@"
$c.Value2 = "2.127"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value2 = "  +3.62%  "   # E25: "  +3.24%  " -> "  +3.62%  "

# Row 26
$c = $ws.Cells.Item(26, 4)   # D26: "20.91" -> "20.88"
$c.NumberFormat = "@"
$c.Value2 = "20.88"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value2 = "  +5.38%  "   # E26: "  +5.43%  " -> "  +5.38%  "

# Row 27
$c = $ws.Cells.Item(27, 4)   # D27: "158.31" -> "158.32"
$c.NumberFormat = "@"
$c.Value2 = "158.32"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value2 = "  +2.45%  "   # E27: "  +2.33%  " -> "  +2.45%  "

# Row 28
$c = $ws.Cells.Item(28, 4)   # D28: "2.446" -> "2.444"
$c.NumberFormat = "@"
$c.Value2 = "2.444"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value2 = "  +18.74%  "   # E28: "  +18.68%  " -> "  +18.74%  "

# Row 29
$ws.Cells.Item(29, 4).Value2 = "2.017.45"   # D29: "2.020.29" -> "2.017.45"
$ws.Cells.Item(29, 5).Value2 = "  +5.19%  "   # E29: "  +5.26%  " -> "  +5.19%  "

# Row 30
$c = $ws.Cells.Item(30, 4)   # D30: "124.05" -> "124.27"
$c.NumberFormat = "@"
$c.Value2 = "124.27"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value2 = "  +3.34%  "   # E30: "  +3.14%  " -> "  +3.34%  "

# Row 31
$c = $ws.Cells.Item(31, 4)   # D31: "1.168" -> "1.169"
$c.NumberFormat = "@"
$c.Value2 = "1.169"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value2 = "  +13.45%  "   # E31: "  +13.04%  " -> "  +13.45%  "

# Row 32
$c = $ws.Cells.Item(32, 4)   # D32: "0.1057" -> "0.1056"
$c.NumberFormat = "@"
$c.Value2 = "0.1056"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = "  +13.92%  "   # E32: "  +13.72%  " -> "  +13.92%  "

# Row 33
$c = $ws.Cells.Item(33, 4)   # D33: "5.791" -> "5.793"
$c.NumberFormat = "@"
$c.Value2 = "5.793"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value2 = "  +8.17%  "   # E33: "  +8.00%  " -> "  +8.17%  "

# Row 34
$c = $ws.Cells.Item(34, 4)   # D34: "3.635" -> "3.632"
$c.NumberFormat = "@"
$c.Value2 = "3.632"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value2 = "  +1.75%  "   # E34: "  +1.82%  " -> "  +1.75%  "

# Row 35
$c = $ws.Cells.Item(35, 4)   # D35: "8.928" -> "8.907"
$c.NumberFormat = "@"
$c.Value2 = "8.907"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = "  +20.02%  "   # E35: "  +20.15%  " -> "  +20.02%  "

# Row 36
$c = $ws.Cells.Item(36, 4)   # D36: "0.02321" -> "0.02322"
$c.NumberFormat = "@"
$c.Value2 = "0.02322"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value2 = "  +6.79%  "   # E36: "  +6.63%  " -> "  +6.79%  "

# Row 37
$c = $ws.Cells.Item(37, 4)   # D37: "0.2168" -> "0.2165"
$c.NumberFormat = "@"
$c.Value2 = "0.2165"
$c.Style = "Normal"

# Row 38
$ws.Cells.Item(38, 5).Value2 = "  +6.99%  "   # E38: "  +7.14%  " -> "  +6.99%  "

# Row 39
$c = $ws.Cells.Item(39, 4)   # D39: "0.06096" -> "0.06098"
$c.NumberFormat = "@"
$c.Value2 = "0.06098"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value2 = "  +3.47%  "   # E39: "  +3.19%  " -> "  +3.47%  "

# Row 40
$c = $ws.Cells.Item(40, 4)   # D40: "0.6418" -> "0.6427"
$c.NumberFormat = "@"
$c.Value2 = "0.6427"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value2 = "  +7.85%  "   # E40: "  +7.72%  " -> "  +7.85%  "

# Row 41
$c = $ws.Cells.Item(41, 4)   # D41: "5.040" -> "5.046"
$c.NumberFormat = "@"
$c.Value2 = "5.046"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value2 = "  +6.88%  "   # E41: "  +6.77%  " -> "  +6.88%  "

# Row 42
$ws.Cells.Item(42, 2).Value2 = "TrustWalletToken"   # B42: "Frax" -> "TrustWalletToken"
$ws.Cells.Item(42, 3).Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"   # C42: "https://coinranking.com/coin/KfWtaeV1W+frax-frax" -> "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Cells.Item(42, 4)   # D42: "0.9989" -> "1.163"
$c.NumberFormat = "@"
$c.Value2 = "1.163"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value2 = "  +4.18%  "   # E42: "  -0.27%  " -> "  +4.18%  "

# Row 43
$ws.Cells.Item(43, 2).Value2 = "Frax"   # B43: "TrustWalletToken" -> "Frax"
$ws.Cells.Item(43, 3).Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"   # C43: "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" -> "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Cells.Item(43, 4)   # D43: "1.159" -> "0.9993"
$c.NumberFormat = "@"
$c.Value2 = "0.9993"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value2 = "  -0.27%  "   # E43: "  +3.82%  " -> "  -0.27%  "

# Row 44
$ws.Cells.Item(44, 5).Value2 = "  -3.80%  "   # E44: "  -3.61%  " -> "  -3.80%  "

# Row 45
$c = $ws.Cells.Item(45, 4)   # D45: "13.41" -> "13.40"
$c.NumberFormat = "@"
$c.Value2 = "13.40"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = "  +5.21%  "   # E45: "  +4.60%  " -> "  +5.21%  "

# Row 46
$c = $ws.Cells.Item(46, 4)   # D46: "0.6007" -> "0.6015"
$c.NumberFormat = "@"
$c.Value2 = "0.6015"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = "  +7.80%  "   # E46: "  +7.65%  " -> "  +7.80%  "

# Row 47
$c = $ws.Cells.Item(47, 4)   # D47: "3.710" -> "3.711"
$c.NumberFormat = "@"
$c.Value2 = "3.711"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value2 = "  +3.85%  "   # E47: "  +3.74%  " -> "  +3.85%  "

# Row 48
$c = $ws.Cells.Item(48, 4)   # D48: "122.38" -> "122.47"
$c.NumberFormat = "@"
$c.Value2 = "122.47"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value2 = "  +2.97%  "   # E48: "  +2.87%  " -> "  +2.97%  "

# Row 49
$ws.Cells.Item(49, 5).Value2 = "  +6.47%  "   # E49: "  +6.40%  " -> "  +6.47%  "

# Row 50
$c = $ws.Cells.Item(50, 4)   # D50: "1.151" -> "1.150"
$c.NumberFormat = "@"
$c.Value2 = "1.150"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = "  +5.29%  "   # E50: "  +5.43%  " -> "  +5.29%  "

# Row 51
$c = $ws.Cells.Item(51, 4)   # D51: "0.06851" -> "0.06854"
$c.NumberFormat = "@"
$c.Value2 = "0.06854"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value2 = "  +3.54%  "   # E51: "  +3.49%  " -> "  +3.54%  "
